# ---------------------------------------------------------------------------
# dunn_matrix_maxsharpe_out_sample_direct_t126.xlsx edit
#   - renames the existing "VaR" sheet to "sharpe_period" and gives it new
#     (recomputed) matrix values
#   - inserts a brand-new "VaR" sheet after "sharpe_period" (end of the
#     workbook) with its own matrix values
#   - refreshes the matrices on "annualised_return", "mean_period_return"
#     (identical data) and "sharpe_annualized" with newly simulated results
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$labels = @(
    "minvar_ports_maxsharpe",
    "minvar_w_cryptos_ahc_maxsharpe",
    "minvar_w_cryptos_kmeans_maxsharpe",
    "minvar_w_cryptos_kshape_maxsharpe",
    "minvar_w_cryptos_random_maxsharpe",
    "rand_ports_maxsharpe",
    "random_w_cryptos_ahc_maxsharpe",
    "random_w_cryptos_kmeans_maxsharpe",
    "random_w_cryptos_kshape_maxsharpe",
    "random_w_cryptos_random_maxsharpe"
)

function New-Matrix {
    param($data)
    $rows = $data.Count
    $cols = $data[0].Count
    $arr = New-Object 'object[,]' $rows, $cols
    for ($r = 0; $r -lt $rows; $r++) {
        for ($c = 0; $c -lt $cols; $c++) {
            $arr[$r, $c] = $data[$r][$c]
        }
    }
    return $arr
}

function Set-DunnMatrix {
    param($ws, $data)
    $m = New-Matrix $data
    $ws.Range("B2:K11").Value = $m
}

function Format-HeaderLabels {
    # Gives a freshly-added, otherwise-empty sheet the same look as the
    # existing dunn-matrix sheets: bold, centered, thin-bordered header row
    # (B1:K1) and left-hand labels (A2:A11), using the workbook's existing
    # 10 share-string labels.
    param($ws)

    for ($i = 0; $i -lt 10; $i++) {
        $col = 2 + $i

        $hcell = $ws.Cells.Item(1, $col)
        $hcell.Value = $labels[$i]
        $hcell.Font.Bold = $true
        $hcell.HorizontalAlignment = -4108   # xlCenter
        $hcell.VerticalAlignment = -4160     # xlTop
        $hcell.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
        $hcell.Borders.Item(8).LineStyle = 1   # xlEdgeTop
        $hcell.Borders.Item(9).LineStyle = 1   # xlEdgeBottom
        $hcell.Borders.Item(10).LineStyle = 1  # xlEdgeRight

        $row = 2 + $i
        $lcell = $ws.Cells.Item($row, 1)
        $lcell.Value = $labels[$i]
        $lcell.Font.Bold = $true
        $lcell.HorizontalAlignment = -4108
        $lcell.VerticalAlignment = -4160
        $lcell.Borders.Item(7).LineStyle = 1
        $lcell.Borders.Item(8).LineStyle = 1
        $lcell.Borders.Item(9).LineStyle = 1
        $lcell.Borders.Item(10).LineStyle = 1
    }
}

# ---------------------------------------------------------------------------
# 1. annualised_return / mean_period_return (identical refreshed matrices)
# ---------------------------------------------------------------------------

$sheet12 = @(
    @(1, 1, 0, 0, 0, 0, 0, 0, 0, 0),
    @(1, 1, 0, 0, 0.012, 0, 0, 0, 0, 0),
    @(0, 0, 1, 1, 1, 0, 0, 0, 0, 0),
    @(0, 0, 1, 1, 1, 0, 0, 0, 0, 0),
    @(0, 0.012, 1, 1, 1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1, 1, 1, 1, 0.1069),
    @(0, 0, 0, 0, 0, 1, 1, 1, 1, 0.0064),
    @(0, 0, 0, 0, 0, 1, 1, 1, 1, 0.0501),
    @(0, 0, 0, 0, 0, 1, 1, 1, 1, 0.9125),
    @(0, 0, 0, 0, 0, 0.1069, 0.0064, 0.0501, 0.9125, 1)
)

Set-DunnMatrix $wb.Worksheets.Item("annualised_return") $sheet12
Set-DunnMatrix $wb.Worksheets.Item("mean_period_return") $sheet12

# ---------------------------------------------------------------------------
# 2. sharpe_annualized
# ---------------------------------------------------------------------------

$sheet3 = @(
    @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 0.2746, 0, 0, 0, 0.341),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0.0105),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0.0107),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0.0002),
    @(0, 0.2746, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0.0003),
    @(0, 0.341, 0.0105, 0.0107, 0.0002, 0, 0, 0, 0.0003, 1)
)

Set-DunnMatrix $wb.Worksheets.Item("sharpe_annualized") $sheet3

# ---------------------------------------------------------------------------
# 3. Old "VaR" sheet becomes "sharpe_period" with freshly computed values,
#    and a brand-new "VaR" sheet is appended with the new VaR results.
# ---------------------------------------------------------------------------

$varSheet = $wb.Worksheets.Item("VaR")
$varSheet.Name = "sharpe_period"

$sheet4 = @(
    @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 1, 1, 1, 1, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0.0002),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0.0002, 1)
)

Set-DunnMatrix $varSheet $sheet4

$newVar = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $varSheet)
$newVar.Name = "VaR"
Format-HeaderLabels $newVar

$sheet5 = @(
    @(1, 0, 0, 0, 0, 0, 0, 0, 0, 0),
    @(0, 1, 0, 0, 0.5101, 0, 0, 0, 0, 0),
    @(0, 0, 1, 1, 0.0492, 0, 0, 0, 0, 0),
    @(0, 0, 1, 1, 0.0487, 0, 0, 0, 0, 0),
    @(0, 0.5101, 0.0492, 0.0487, 1, 0, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 1, 0, 0, 0, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 1, 1, 1, 0),
    @(0, 0, 0, 0, 0, 0, 0, 0, 0, 1)
)

Set-DunnMatrix $newVar $sheet5
